$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "#PBM:DWP:Request:Check user roles and company association for 'On behalf of' requests.:Educate"
$ws.Range("E3").Value = "#PBM:DWP:Data:Mismatch in ticket status between user profile and Service Desk.:NA"
$ws.Range("E6").Value = "#PBM:DWP:Integration:Guide on automating user access via API.:Educate"
$ws.Range("E8").Value = "#PBM:DWP:Request:Investigate request cancellation delay in DWP.:R&D"
$ws.Range("E15").Value = "#PBM:DWP:Multitenancy:Enhanced catalog setup for sub tenant on production.:Customization"
$ws.Range("E18").Value = "#PBM:DWP:Broadcast:Broadcast not removed due to sync issue; needs manual update.:NA"
$ws.Range("E20").Value = "#PBM:DWP:Customization:Usernames in comments can be customized via settings.:Educate"
